$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "PhD in Neuroscience  - \textbf{\textit{Summa Cum Laude}}"
$ws.Range("A4").Value = "Psychological Research Methods (Evolutionary Psychology) MSc - \textbf{\textit{Distinction}}"

$ws.Range("A2").Select()
